# Applies the commit's edit to the Magic 8 Ball document:
#  1. Adds an (empty) "_GoBack" bookmark at the end of the paragraph
#     that ends with "...you don't do things for free.  " (Word stamps
#     this bookmark at the location of the most recent edit).
#  2. Removes the final paragraph's picture (Picture 3 / the screenshot
#     that had the page-break + the old "_GoBack" bookmark), leaving
#     that paragraph empty - this is where the "_GoBack" bookmark used
#     to live before the edit above moved it.

$d = $word.ActiveDocument

# --- 1. Move/create the "_GoBack" bookmark onto the community-service paragraph ---

$target = $d.Content
$found = $target.Find.Execute("do things for free.  ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # $target now covers the matched text; collapsing straight to the very
    # end of the paragraph (a zero-length point that sits on the paragraph
    # boundary) mis-resolves to the *next* paragraph's start, so instead we
    # briefly park a marker character there, bookmark that exact offset,
    # then remove the marker again - leaving the bookmark correctly inside
    # the original paragraph, right after the last run.
    $target.Collapse(0)
    $target.InsertAfter("~MARK~")

    $markStart = $target.Start
    $markEnd = $markStart + 6

    $bmRange = $d.Range($markStart, $markStart)
    $d.Bookmarks.Add("_GoBack", $bmRange)

    $d.Range($markStart, $markEnd).Delete()
}

# --- 2. Strip the last paragraph back down to an empty paragraph ---

$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)

# Remove any inline pictures that live in this paragraph (walk backwards
# since deleting shifts indices).
$shapes = $d.InlineShapes
for ($i = $shapes.Count; $i -ge 1; $i--) {
    $shape = $shapes.Item($i)
    $shapeRange = $shape.Range
    if ($shapeRange.Start -ge $lastPara.Range.Start -and $shapeRange.Start -lt $lastPara.Range.End) {
        $shape.Delete()
    }
}

# Re-fetch (index/range may have shifted) now that the picture run is gone.
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)

# Clear paragraph formatting (e.g. centered alignment) back to default.
$lastPara.Format.Alignment = 0

# Defensive cleanup: if a "_GoBack" bookmark is still anchored in this
# paragraph (e.g. step 1 above didn't find/move it for some reason), drop
# it too so the paragraph ends up completely empty.
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    if ($bm -ne $null) {
        $bmStart = $bm.Start
        if ($bmStart -ge $lastPara.Range.Start -and $bmStart -le $lastPara.Range.End) {
            $bm.Delete()
        }
    }
} catch {
    # No "_GoBack" bookmark at all - nothing to clean up.
}

Write-Output "done"
